$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels (bold, matching the existing Min/Max, Q1/Median, Q3/IQR pairs)
$ws.Range("D15").Value = "Std"
$ws.Range("E15").Value = "Relative std"
$ws.Range("D15:E15").Font.Bold = $true

# New statistics: standard deviation and relative standard deviation (%)
$ws.Range("D16").Formula = "=STDEV(B2:B31)"
$ws.Range("E16").Formula = "=(D16/E4)*100"

# Update the sheet's active selection
$ws.Range("H21").Select() | Out-Null
